$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SheetName1")
$ws.Range("T2").Value = 73959
